# Update the "按登记注册类型分建筑业企业税金总额" sheet:
#  - drop the oldest years (2005-2009), shifting the remaining rows up
#  - append a new row for 2021 at the bottom
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows for 2005年-2009年 (original rows 2-6); remaining rows shift up,
# so 2010年 (was row 7) becomes row 2, ..., 2020年 (was row 16) becomes row 11.
$ws.Range("A2:I6").Delete() | Out-Null

# Copy the formatting of the preceding year-label cell onto the new row
# so the appended row matches the existing look (bold, bordered, centered).
$ws.Range("A11").Copy($ws.Range("A12")) | Out-Null

# Append the new 2021年 row of data.
$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 6997.2087
$ws.Range("C12").Value = 712.64289
$ws.Range("D12").Value = 28.80583
$ws.Range("E12").Value = 16.57063
$ws.Range("F12").Value = 7052.36398
$ws.Range("G12").Value = 26.34945
$ws.Range("H12").Value = 4.42923
$ws.Range("I12").Value = 144.17168
